$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Swap B5 and C5 values (Kentucky <-> Georgia)
$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2
$ws.Range("B5").Value = $c5
$ws.Range("C5").Value = $b5

# Set best-fit column widths for columns B and C (matches Excel's
# "double-click the column border" autofit result for these columns)
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 13.5

# Update selection to C6
$ws.Range("C6").Select()
